# Undo Jason's overwrite of class materials starter code/slides.
#
# 1) Restore the "Today's Attendance password" textbox on slide 1: the
#    placeholder word "timers" becomes a blank-style "_________".
# 2) Roll the cached text of the auto "today" date field back from
#    9/3/2023 to 1/23/2023 everywhere it is stamped (slide master and
#    every slide layout).

$p = $ppt.ActivePresentation

# --- 1) Slide 1 text box: "timers" -> "_________" ------------------------
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.HasTextFrame -eq -1) {
        $tr = $shp.TextFrame.TextRange
        for ($pi = 1; $pi -le $tr.Paragraphs().Count; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            if ($para.Text -eq "timers") {
                $para.Text = "_________"
            }
        }
    }
}

# --- 2) Fix up the cached "today" date stamped on masters/layouts --------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "9/3/2023") {
                $tr.Text = "1/23/2023"
            }
        }
    }
}

Update-DatePlaceholders $p.SlideMaster.Shapes
Update-DatePlaceholders $p.NotesMaster.Shapes
Update-DatePlaceholders $p.HandoutMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes
}
